$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- original row 16
$ws.Range("D2").Value = 44204
$ws.Range("K2").Value = 500
$ws.Range("L2").Value = 550
$ws.Range("M2").Value = 525
$ws.Range("P2").Value = 525

# Row 3 <- original row 7
$ws.Range("D3").Value = 44575
$ws.Range("J3").Value = 1300
$ws.Range("K3").Value = 500
$ws.Range("L3").Value = 550
$ws.Range("M3").Value = 525
$ws.Range("O3").Value = 'Región de O''Higgins'
$ws.Range("P3").Value = 525

# Row 4 <- original row 2
$ws.Range("D4").Value = 44243
$ws.Range("I4").Value = '2a nueva(o)'
$ws.Range("J4").Value = 1600
$ws.Range("K4").Value = 450
$ws.Range("L4").Value = 480
$ws.Range("M4").Value = 465
$ws.Range("O4").Value = 'Región del Maule'
$ws.Range("P4").Value = 465

# Row 5 <- original row 24
$ws.Range("D5").Value = 44490
$ws.Range("I5").Value = '1a nueva(o)'
$ws.Range("J5").Value = 1200
$ws.Range("K5").Value = 450
$ws.Range("L5").Value = 480
$ws.Range("M5").Value = 465
$ws.Range("P5").Value = 465

# Row 6 <- original row 15
$ws.Range("D6").Value = 44201
$ws.Range("J6").Value = 1360
$ws.Range("K6").Value = 730
$ws.Range("L6").Value = 750
$ws.Range("M6").Value = 740
$ws.Range("P6").Value = 740

# Row 7 <- original row 22
$ws.Range("D7").Value = 44301
$ws.Range("I7").Value = '2a nueva(o)'
$ws.Range("J7").Value = 1200
$ws.Range("K7").Value = 400
$ws.Range("L7").Value = 430
$ws.Range("M7").Value = 415
$ws.Range("O7").Value = 'Provincia de Melipilla'
$ws.Range("P7").Value = 415

# Row 8 <- original row 18
$ws.Range("D8").Value = 44179
$ws.Range("J8").Value = 1000
$ws.Range("K8").Value = 1350
$ws.Range("L8").Value = 1400
$ws.Range("M8").Value = 1375
$ws.Range("P8").Value = 1375

# Row 9 <- original row 20
$ws.Range("D9").Value = 44483
$ws.Range("I9").Value = '1a nueva(o)'
$ws.Range("J9").Value = 1300
$ws.Range("K9").Value = 550
$ws.Range("L9").Value = 580
$ws.Range("M9").Value = 565
$ws.Range("O9").Value = 'Perú'
$ws.Range("P9").Value = 565

# Row 10 <- original row 6
$ws.Range("D10").Value = 44469
$ws.Range("J10").Value = 1200
$ws.Range("K10").Value = 600
$ws.Range("L10").Value = 650
$ws.Range("M10").Value = 625
$ws.Range("P10").Value = 625

# Row 11 <- original row 17
$ws.Range("D11").Value = 44168
$ws.Range("I11").Value = '1a nueva(o)'
$ws.Range("J11").Value = 1200
$ws.Range("K11").Value = 1500
$ws.Range("L11").Value = 1700
$ws.Range("M11").Value = 1600
$ws.Range("P11").Value = 1600

# Row 12 <- original row 13
$ws.Range("D12").Value = 44349
$ws.Range("H12").Value = 'Pachia'
$ws.Range("K12").Value = 730
$ws.Range("L12").Value = 750
$ws.Range("M12").Value = 740
$ws.Range("P12").Value = 740

# Row 13 <- original row 10
$ws.Range("D13").Value = 44231
$ws.Range("H13").Value = 'Camote'
$ws.Range("J13").Value = 1300
$ws.Range("K13").Value = 450
$ws.Range("L13").Value = 480
$ws.Range("M13").Value = 465
$ws.Range("P13").Value = 465

# Row 14 <- original row 21
$ws.Range("D14").Value = 44476
$ws.Range("J14").Value = 1200
$ws.Range("K14").Value = 480
$ws.Range("L14").Value = 500
$ws.Range("M14").Value = 490
$ws.Range("P14").Value = 490

# Row 15 <- original row 4
$ws.Range("D15").Value = 44547
$ws.Range("I15").Value = '1a (cosecha)'
$ws.Range("J15").Value = 800
$ws.Range("K15").Value = 600
$ws.Range("L15").Value = 650
$ws.Range("M15").Value = 625
$ws.Range("P15").Value = 625

# Row 16 <- original row 5
$ws.Range("D16").Value = 44547
$ws.Range("J16").Value = 300
$ws.Range("O16").Value = 'Perú'

# Row 17 <- original row 3
$ws.Range("D17").Value = 44238
$ws.Range("J17").Value = 1250
$ws.Range("K17").Value = 430
$ws.Range("L17").Value = 450
$ws.Range("M17").Value = 440
$ws.Range("P17").Value = 440

# Row 18 <- original row 11
$ws.Range("D18").Value = 44530
$ws.Range("I18").Value = '2a nueva(o)'
$ws.Range("J18").Value = 900
$ws.Range("K18").Value = 480
$ws.Range("L18").Value = 500
$ws.Range("M18").Value = 490
$ws.Range("P18").Value = 490

# Row 19 <- original row 23
$ws.Range("D19").Value = 44322
$ws.Range("I19").Value = '1a (cosecha)'
$ws.Range("J19").Value = 1200
$ws.Range("K19").Value = 350
$ws.Range("L19").Value = 400
$ws.Range("M19").Value = 375
$ws.Range("O19").Value = 'Región del Maule'
$ws.Range("P19").Value = 375

# Row 20 <- original row 19
$ws.Range("D20").Value = 44211
$ws.Range("J20").Value = 1600
$ws.Range("K20").Value = 500
$ws.Range("L20").Value = 550
$ws.Range("M20").Value = 525
$ws.Range("O20").Value = 'Región de O''Higgins'
$ws.Range("P20").Value = 525

# Row 21 <- original row 12
$ws.Range("D21").Value = 44175
$ws.Range("K21").Value = 1500
$ws.Range("L21").Value = 1550
$ws.Range("M21").Value = 1525
$ws.Range("P21").Value = 1525

# Row 22 <- original row 9
$ws.Range("D22").Value = 44399
$ws.Range("I22").Value = '1a (guarda)'
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 450
$ws.Range("L22").Value = 480
$ws.Range("M22").Value = 465
$ws.Range("P22").Value = 465

# Row 23 <- original row 14
$ws.Range("D23").Value = 44453
$ws.Range("I23").Value = '1a nueva(o)'
$ws.Range("J23").Value = 800
$ws.Range("K23").Value = 630
$ws.Range("L23").Value = 650
$ws.Range("M23").Value = 640
$ws.Range("O23").Value = 'Perú'
$ws.Range("P23").Value = 640

# Row 24 <- original row 8
$ws.Range("D24").Value = 44428
$ws.Range("J24").Value = 600
$ws.Range("K24").Value = 580
$ws.Range("L24").Value = 600
$ws.Range("M24").Value = 590
$ws.Range("P24").Value = 590
